# Updated symbol list refresh: Price (column D) and Volume(1h) (column E)
# values for the coin rows on the "cryptos" worksheet, matching the
# GitHub Actions data-refresh commit.
#
# These cells are stored as text (not numbers/percentages) in the
# workbook, so we force the cell NumberFormat to "@" (Text) before
# assigning the new values -- this prevents Excel from silently
# re-interpreting strings like "292.70" or "0.68%" as numeric/percent
# values and mangling their literal text (e.g. dropping trailing
# zeros or converting "0.68%" into 0.0068). After the values are set
# we restore the cell style back to "Normal" so no stray formatting
# is left behind on the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetCells = @(
    "D2","E2",
    "D3","E3",
    "D4","E4",
    "D5","E5",
    "D6","E6",
    "D7","E7",
    "D8","E8",
    "D9","E9",
    "D10","E10",
    "D11","E11",
    "D12","E12",
    "D13","E13",
    "D14","E14",
    "D15","E15",
    "D16","E16",
    "D17","E17",
    "D18","E18",
    "E19",
    "D20","E20",
    "D21","E21",
    "D22","E22",
    "D23","E23",
    "D24","E24",
    "E25",
    "D26","E26",
    "D27",
    "D39","E39",
    "D40","E40",
    "D41","E41",
    "D42","E42",
    "D43","E43",
    "D44","E44",
    "D45","E45",
    "D46","E46",
    "D47","E47",
    "D48","E48",
    "D49","E49",
    "D50","E50",
    "D51","E51"
)

foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value  = "292.70"
$ws.Range("E2").Value  = "0.68%"
$ws.Range("D3").Value  = "31.01"
$ws.Range("E3").Value  = "0.53%"
$ws.Range("D4").Value  = "4.925"
$ws.Range("E4").Value  = "-0.38%"
$ws.Range("D5").Value  = "0.07397"
$ws.Range("E5").Value  = "2.77%"
$ws.Range("D6").Value  = "2.193"
$ws.Range("E6").Value  = "21.95%"
$ws.Range("D7").Value  = "7.732"
$ws.Range("E7").Value  = "0.88%"
$ws.Range("D8").Value  = "3.750"
$ws.Range("E8").Value  = "0.10%"
$ws.Range("D9").Value  = "0.9141"
$ws.Range("E9").Value  = "2.08%"
$ws.Range("D10").Value = "0.08952"
$ws.Range("E10").Value = "15.88%"
$ws.Range("D11").Value = "0.1703"
$ws.Range("E11").Value = "3.01%"
$ws.Range("D12").Value = "0.08298"
$ws.Range("E12").Value = "3.38%"
$ws.Range("D13").Value = "0.03125"
$ws.Range("E13").Value = "2.01%"
$ws.Range("D14").Value = "0.09983"
$ws.Range("E14").Value = "-0.45%"
$ws.Range("D15").Value = "0.001505"
$ws.Range("E15").Value = "1.04%"
$ws.Range("D16").Value = "0.005849"
$ws.Range("E16").Value = "1.38%"
$ws.Range("D17").Value = "3.499"
$ws.Range("E17").Value = "0.79%"
$ws.Range("D18").Value = "2.167"
$ws.Range("E18").Value = "4.17%"
$ws.Range("E19").Value = "1.47%"
$ws.Range("D20").Value = "0.1297"
$ws.Range("E20").Value = "1.88%"
$ws.Range("D21").Value = "3.988"
$ws.Range("E21").Value = "-1.40%"
$ws.Range("D22").Value = "0.2188"
$ws.Range("E22").Value = "9.45%"
$ws.Range("D23").Value = "0.04558"
$ws.Range("E23").Value = "0.76%"
$ws.Range("D24").Value = "0.001215"
$ws.Range("E24").Value = "0.23%"
$ws.Range("E25").Value = "14.36%"
$ws.Range("D26").Value = "0.0001303"
$ws.Range("E26").Value = "4.25%"
$ws.Range("D27").Value = "0.0003399"
$ws.Range("D39").Value = "0.01600"
$ws.Range("E39").Value = "-0.20%"
$ws.Range("D40").Value = "0.04482"
$ws.Range("E40").Value = "2.18%"
$ws.Range("D41").Value = "0.007334"
$ws.Range("E41").Value = "-0.18%"
$ws.Range("D42").Value = "0.009668"
$ws.Range("E42").Value = "26.32%"
$ws.Range("D43").Value = "0.1327"
$ws.Range("E43").Value = "1.56%"
$ws.Range("D44").Value = "0.002325"
$ws.Range("E44").Value = "12.34%"
$ws.Range("D45").Value = "0.009141"
$ws.Range("E45").Value = "-0.72%"
$ws.Range("D46").Value = "0.00006095"
$ws.Range("E46").Value = "3.16%"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "0.13%"
$ws.Range("D48").Value = "2.206"
$ws.Range("E48").Value = "-1.76%"
$ws.Range("D49").Value = "0.002002"
$ws.Range("E49").Value = "-33.24%"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").Value = "0.13%"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").Value = "0.13%"

foreach ($addr in $targetCells) {
    $ws.Range($addr).Style = "Normal"
}
